$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.235.46"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.155.65"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.145.43"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "3.684.75"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "64.042.17"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "3.154.84"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.61%  "
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.46%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.113"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.01%  "
$ws.Range("D35").Value = "0.0₃0869"
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.27%  "
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "461.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.89%  "
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").Value = "2.908.41"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.39%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.57%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.56%  "
